$d = $word.ActiveDocument

# Locate the title run's text ("Analysis of the virality of Reddit posts")
# and grab its whole containing paragraph (including the paragraph mark)
# so we can swap it for a version split across three runs.
$hit = $d.Content
$found = $hit.Find.Execute("Analysis of the virality of Reddit posts")
if (-not $found) {
    throw "Could not find the title text to edit."
}
$titlePara = $hit.Paragraphs(1)
$paraRange = $titlePara.Range

# Rebuild the paragraph as literal OOXML so the new text lands in three
# separate <w:r> runs (matching how the authored edit split "virality"
# -> " popularity" across new runs) instead of being coalesced back into
# a single run the way plain text/Find-Replace edits would be.
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14 = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$rPr = '<w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'

$xml = '<w:p ' + $w + ' ' + $w14 + ' w14:paraId="2D957B26" w14:textId="720064CB" w:rsidR="00162C06" w:rsidRDefault="00162C06" w:rsidP="00162C06">' +
       '<w:pPr><w:spacing w:before="360" w:after="120"/><w:outlineLvl w:val="1"/>' + $rPr + '</w:pPr>' +
       '<w:r w:rsidRPr="00162C06">' + $rPr + '<w:t>Analysis of the</w:t></w:r>' +
       '<w:r>' + $rPr + '<w:t xml:space="preserve"> popularity</w:t></w:r>' +
       '<w:r>' + $rPr + '<w:t xml:space="preserve"> of Reddit posts</w:t></w:r>' +
       '</w:p>'

$paraRange.InsertXML($xml)
